# Apply weekly fruit/vegetable price update: rows are re-associated with
# a different date's figures (rows permuted by date), leaving all other
# columns (Mercado, Región, Codreg, Categoría, Variedad, Calidad, Unidad
# de comercialización, Origen, Kg o Unidades, Clasificación) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
# Values below are the Excel date serial numbers / numbers as they must
# appear in the saved workbook.

$updates = @{
    2 = @{ D = 44389; J = 120; K = 12000; L = 13000; M = 12500; P = 962 }
    3 = @{ D = 44229; J = 120; K = 44000; L = 45000; M = 44500; P = 3423 }
    4 = @{ D = 44159; J = 100; K = 23000; L = 24000; M = 23500; P = 1808 }
    5 = @{ D = 44320; J = 160; K = 19000; L = 20000; M = 19500; P = 1500 }
    6 = @{ D = 44406; J = 160; K = 17000; L = 18000; M = 17500; P = 1346 }
    7 = @{ D = 44397; J = 140; K = 12500; L = 13000; M = 12750; P = 981 }
    8 = @{ D = 44379; J = 120; K = 12000; L = 13000; M = 12667; P = 974 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
